$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.269143104553223
$ws.Range("B1").Value = 1.286441206932068
$ws.Range("C1").Value = 1.386765599250793
$ws.Range("D1").Value = 1.997619867324829
$ws.Range("E1").Value = 4.086225032806396
